$wb = $excel.ActiveWorkbook

# Sheet "Bico": clear the "Obs_relatorio" column (H2:H11) validation messages
$bico = $wb.Worksheets.Item("Bico")
for ($r = 2; $r -le 11; $r++) {
    $bico.Cells.Item($r, 8).Value = ""
}

# Sheet "Tanque": update the "Obs_relatorio" column (F2:F5) with divergence messages
$tanque = $wb.Worksheets.Item("Tanque")
$tanque.Cells.Item(2, 6).Value = "Divergência entre o SPED(11336,50) e o relatório(9134215,00)!"
$tanque.Cells.Item(3, 6).Value = "Divergência entre o SPED(11336,50) e o relatório(3283312,00)!"
$tanque.Cells.Item(4, 6).Value = "Divergência entre o SPED(11336,50) e o relatório(6488804,00)!"
$tanque.Cells.Item(5, 6).Value = "Divergência entre o SPED(11336,50) e o relatório(11336496,00)!"
